$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage (so numeric-looking
# strings like "555.49" are NOT silently coerced into numbers by Excel's
# usual "looks like a number" auto-detection) and without leaving the
# cell's visible style/format changed afterwards.
function Set-TextValue($rangeAddr, $val) {
    $c = $ws.Range($rangeAddr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue "D2" "64.589.71"
Set-TextValue "E2" "  -0.42%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.361.71"
Set-TextValue "E3" "  -2.19%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "555.49"
Set-TextValue "E5" "  -3.06%  "

# Row 6 - Solana
Set-TextValue "D6" "175.82"
Set-TextValue "E6" "  +0.54%  "

# Row 7 - XRP
Set-TextValue "E7" "  -0.39%  "

# Row 8 - LidoStakedEther
Set-TextValue "D8" "3.351.78"
Set-TextValue "E8" "  -2.26%  "

# Row 9 - USDC
Set-TextValue "E9" "  +0.01%  "

# Row 10 - swap Dogecoin -> Cardano
Set-TextValue "B10" "Cardano"
Set-TextValue "C10" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D10" "0.630"
Set-TextValue "E10" "  +0.72%  "

# Row 11 - swap Cardano -> Dogecoin
Set-TextValue "B11" "Dogecoin"
Set-TextValue "C11" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D11" "0.164"
Set-TextValue "E11" "  +3.53%  "

# Row 12 - Avalanche
Set-TextValue "D12" "54.46"
Set-TextValue "E12" "  -1.22%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000274"
Set-TextValue "E13" "  +0.61%  "

# Row 14 - Polkadot
Set-TextValue "E14" "  -0.50%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "3.895.97"
Set-TextValue "E15" "  -2.16%  "

# Row 16 - Chainlink
Set-TextValue "E16" "  +1.72%  "

# Row 17 - TRON
Set-TextValue "E17" "  -2.04%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "3.358.29"
Set-TextValue "E18" "  -2.61%  "

# Row 19 - Uniswap
Set-TextValue "E19" "  -0.22%  "

# Row 20 - WrappedBTC
Set-TextValue "D20" "64.514.37"
Set-TextValue "E20" "  -0.55%  "

# Row 21 - Polygon
Set-TextValue "D21" "0.987"
Set-TextValue "E21" "  -0.42%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "463.59"
Set-TextValue "E22" "  +13.64%  "

# Row 23 - Toncoin
Set-TextValue "D23" "4.79"
Set-TextValue "E23" "  +10.46%  "

# Row 24 - PancakeSwap
Set-TextValue "E24" "  -2.68%  "

# Row 25 - Litecoin
Set-TextValue "D25" "86.17"
Set-TextValue "E25" "  +3.21%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "13.38"
Set-TextValue "E26" "  +1.27%  "

# Row 27 - RenderToken
Set-TextValue "D27" "10.94"
Set-TextValue "E27" "  +0.95%  "

# Row 28 - ImmutableX
Set-TextValue "E28" "  +1.67%  "

# Row 29 - Filecoin
Set-TextValue "D29" "8.79"
Set-TextValue "E29" "  -1.74%  "

# Row 30 - EthereumClassic
Set-TextValue "D30" "30.13"
Set-TextValue "E30" "  +0.82%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "6.65"
Set-TextValue "E31" "  -0.44%  "

# Row 32 - Cosmos
Set-TextValue "E32" "  -0.46%  "

# Row 33 - Bittensor
Set-TextValue "D33" "582.03"
Set-TextValue "E33" "  -0.73%  "

# Row 34 - Hedera
Set-TextValue "E34" "  -0.05%  "

# Row 35 - OKB
Set-TextValue "D35" "58.87"
Set-TextValue "E35" "  -0.41%  "

# Row 36 - Dai
Set-TextValue "E36" "  +0.06%  "

# Row 37 - Kaspa
Set-TextValue "E37" "  -8.27%  "

# Row 38 - Stacks
Set-TextValue "E38" "  -1.02%  "

# Row 39 - PEPE
Set-TextValue "D39" "0.0₃0759"
Set-TextValue "E39" "  -1.61%  "

# Row 40 - InjectiveProtocol
Set-TextValue "E40" "  -1.41%  "

# Row 41 - TheGraph
Set-TextValue "E41" "  -0.23%  "

# Row 42 - Maker
Set-TextValue "D42" "3.103.25"
Set-TextValue "E42" "  -2.49%  "

# Row 43 - FirstDigitalUSD
Set-TextValue "D43" "0.998"
Set-TextValue "E43" "  +0.01%  "

# Row 44 - Fetch.AI
Set-TextValue "E44" "  +1.17%  "

# Row 45 - ThetaToken
Set-TextValue "D45" "2.80"
Set-TextValue "E45" "  -4.15%  "

# Row 46 - VeChain
Set-TextValue "D46" "0.0412"
Set-TextValue "E46" "  +0.31%  "

# Row 47 - ApeXProtocol
Set-TextValue "E47" "  -0.58%  "

# Row 48 - Stellar
Set-TextValue "E48" "  +0.33%  "

# Row 49 - WEMIXToken
Set-TextValue "E49" "  -2.27%  "

# Row 50 - THORChain
Set-TextValue "D50" "8.39"
Set-TextValue "E50" "  -0.61%  "

# Row 51 - Monero
Set-TextValue "D51" "135.42"
Set-TextValue "E51" "  -0.64%  "
